# Apply the "solicitudes_entrantes" update:
#  - row 2 becomes a brand new request (Ariel Gomez Cifuentes / Drone Parrot Bebop 2)
#  - rows 3-9 shift up (their ID de detalle / other fields are re-sequenced)
#  - rows 10-12 (the old tail of the table) are removed
#  - columns E and H get narrower
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: write a numeric-looking value as TEXT without disturbing the
# cell's existing style (a plain .Value = "199029932" would be auto-coerced
# to a number). We stage it as a text formula, then collapse the formula to
# its literal value with a values-only paste, which keeps style/format as-is.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}
$excelApp = $excel

# ---- Row 2 : new request ----
$ws.Cells.Item(2,1).Value = 270
$ws.Cells.Item(2,2).Value = 849
Set-TextValue 2 3 "199029932"
$ws.Cells.Item(2,4).Value = "Ariel Gomez Cifuentes"
$ws.Cells.Item(2,5).Value = "Drone Parrot Bebop 2"
$ws.Cells.Item(2,6).Value = "ABC432"
$ws.Cells.Item(2,7).Value = 44112.0656712963
$ws.Cells.Item(2,8).ClearContents()

# ---- Row 3 ----
$ws.Cells.Item(3,1).Value = 265
$ws.Cells.Item(3,2).Value = 832
Set-TextValue 3 3 "188639089"
$ws.Cells.Item(3,4).Value = "Lino Cisternas"
$ws.Cells.Item(3,5).Value = "Router sony PS5`t Sony PS5"
$ws.Cells.Item(3,6).Value = "SONYPS"
$ws.Cells.Item(3,7).Value = 44111.03188657408
$ws.Cells.Item(3,8).ClearContents()

# ---- Row 4 ----
$ws.Cells.Item(4,1).Value = 265
$ws.Cells.Item(4,2).Value = 834
Set-TextValue 4 3 "188639089"
$ws.Cells.Item(4,4).Value = "Lino Cisternas"
$ws.Cells.Item(4,5).Value = "Router sony PS5`t Sony PS5"
$ws.Cells.Item(4,6).Value = "SONYPS"
$ws.Cells.Item(4,7).Value = 44111.03188657408
$ws.Cells.Item(4,8).ClearContents()

# ---- Row 5 : only "ID de detalle" changes ----
$ws.Cells.Item(5,2).Value = 831

# ---- Row 6 : only "ID de detalle" changes ----
$ws.Cells.Item(6,2).Value = 833

# ---- Row 7 : only "ID de detalle" changes ----
$ws.Cells.Item(7,2).Value = 830

# ---- Row 8 ----
$ws.Cells.Item(8,1).Value = 264
$ws.Cells.Item(8,2).Value = 828
$ws.Cells.Item(8,7).Value = 44111.02668981482

# ---- Row 9 ----
$ws.Cells.Item(9,1).Value = 263
$ws.Cells.Item(9,2).Value = 826
$ws.Cells.Item(9,3).Value = "19889608K"
$ws.Cells.Item(9,4).Value = "Sebastián Ignacio Toro Severino"
$ws.Cells.Item(9,5).Value = "Drone Parrot Bebop 2"
$ws.Cells.Item(9,6).Value = "ABC432"
$ws.Cells.Item(9,7).Value = 44110.98806712963

$excelApp.CutCopyMode = $false

# ---- remove the old rows 10-12 entirely ----
$ws.Rows("10:12").Delete()

# ---- narrower columns E (30 -> 25) and H (28 -> 19) ----
# ColumnWidth (chars) reported by Excel is the OOXML <col width> minus
# 5/6ths of a character, so subtract 0.8333333 to land exactly on target.
$ws.Columns("E").ColumnWidth = 24.1666667
$ws.Columns("H").ColumnWidth = 18.1666667
